$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 475.58334
$ws.Range("J2").Value = 488.8
$ws.Range("L2").Value = 488.8
$ws.Range("N2").Value = -714.8

$ws.Range("H32").Value = 795.53845
$ws.Range("I32").Value = 941.3333
$ws.Range("J32").Value = 670.5714
$ws.Range("K32").Value = 941.3333
$ws.Range("L32").Value = 670.5714
$ws.Range("M32").Value = -615.3333
$ws.Range("N32").Value = -1322.5714

$ws.Range("H55").Value = 298.22223
$ws.Range("I55").Value = 138.4
$ws.Range("J55").Value = 498
$ws.Range("K55").Value = 138.4
$ws.Range("L55").Value = 498
$ws.Range("M55").Value = 75.59999999999999
$ws.Range("N55").Value = -926

$ws.Range("H112").Value = 2251.1177
$ws.Range("I112").Value = 1607.8
$ws.Range("J112").Value = 2519.1667
$ws.Range("K112").Value = 4823.4
$ws.Range("L112").Value = 7557.500100000001
$ws.Range("M112").Value = -3715.4
$ws.Range("N112").Value = -9773.500100000001

$ws.Range("H129").Value = 1833.6666
$ws.Range("I129").Value = 1497.4
$ws.Range("K129").Value = 4492.200000000001
$ws.Range("M129").Value = 507.7999999999993

$ws.Range("H132").Value = 2086589.8
$ws.Range("I132").Value = 3100.25
$ws.Range("J132").Value = 25004976
$ws.Range("K132").Value = 9300.75
$ws.Range("L132").Value = 75014928
$ws.Range("M132").Value = -6770.75
$ws.Range("N132").Value = -75019988

$ws.Range("H137").Value = 11161.096
$ws.Range("I137").Value = 13768.625
$ws.Range("K137").Value = 41305.875
$ws.Range("M137").Value = -38755.875

$ws.Range("H138").Value = 211322.78
$ws.Range("I138").Value = 388096.2
$ws.Range("J138").Value = 4071.2068
$ws.Range("K138").Value = 1164288.6
$ws.Range("L138").Value = 12213.6204
$ws.Range("M138").Value = -1159148.6
$ws.Range("N138").Value = -22493.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5115.817
$ws.Range("I32").Value = 5091.5146
$ws.Range("J32").Value = 5666.6665
$ws.Range("K32").Value = 5091.5146
$ws.Range("L32").Value = 5666.6665
$ws.Range("M32").Value = -4804.5146
$ws.Range("N32").Value = -6240.6665

$ws.Range("H60").Value = 3748.5
$ws.Range("I60").Value = 3748.5
$ws.Range("K60").Value = 3748.5
$ws.Range("M60").Value = -3015.5

$ws.Range("H110").Value = 1884
$ws.Range("I110").Value = 1511.909
$ws.Range("K110").Value = 1511.909
$ws.Range("M110").Value = 533.0909999999999

$ws.Range("H122").Value = 1192942.2
$ws.Range("I122").Value = 3359.8076
$ws.Range("J122").Value = 4004682.5
$ws.Range("K122").Value = 10079.4228
$ws.Range("L122").Value = 12014047.5
$ws.Range("M122").Value = -7629.4228
$ws.Range("N122").Value = -12018947.5

$ws.Range("H132").Value = 2062.6123
$ws.Range("J132").Value = 3017.923
$ws.Range("L132").Value = 9053.769
$ws.Range("N132").Value = -14113.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 13661.964
$ws.Range("I99").Value = 15749.048
$ws.Range("K99").Value = 15749.048
$ws.Range("M99").Value = -14251.048

$ws.Range("H105").Value = 104829.4
$ws.Range("I105").Value = 167882.33
$ws.Range("K105").Value = 167882.33
$ws.Range("M105").Value = -166135.33

$ws.Range("H134").Value = 5201.4688
$ws.Range("I134").Value = 6674.478
$ws.Range("J134").Value = 1437.1111
$ws.Range("K134").Value = 20023.434
$ws.Range("L134").Value = 4311.3333
$ws.Range("M134").Value = -17488.434
$ws.Range("N134").Value = -9381.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6768.72
$ws.Range("I31").Value = 6768.72
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6768.72
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6473.72

$ws.Range("H34").Value = 6768.72
$ws.Range("I34").Value = 6768.72
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6768.72
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6566.72

$ws.Range("H58").Value = 3179.625
$ws.Range("I58").Value = 2951
$ws.Range("J58").Value = 3560.6667
$ws.Range("K58").Value = 2951
$ws.Range("L58").Value = 3560.6667
$ws.Range("M58").Value = -2748
$ws.Range("N58").Value = -3966.6667

$ws.Range("H99").Value = 12504770
$ws.Range("I99").Value = 20836700
$ws.Range("K99").Value = 20836700
$ws.Range("M99").Value = -20835202

$ws.Range("H122").Value = 10911.692
$ws.Range("I122").Value = 18540.285
$ws.Range("K122").Value = 55620.855
$ws.Range("M122").Value = -53170.855

$ws.Range("H126").Value = 12504770
$ws.Range("I126").Value = 20836700
$ws.Range("K126").Value = 62510100
$ws.Range("M126").Value = -62507630

$ws.Range("H132").Value = 1471.52
$ws.Range("I132").Value = 1367.6818
$ws.Range("J132").Value = 2233
$ws.Range("K132").Value = 4103.0454
$ws.Range("L132").Value = 6699
$ws.Range("M132").Value = -1573.0454
$ws.Range("N132").Value = -11759

$ws.Range("H134").Value = 2391.0938
$ws.Range("I134").Value = 2490.75
$ws.Range("K134").Value = 7472.25
$ws.Range("M134").Value = -4937.25

$ws.Range("H136").Value = 3179.625
$ws.Range("I136").Value = 2951
$ws.Range("J136").Value = 3560.6667
$ws.Range("K136").Value = 8853
$ws.Range("L136").Value = 10682.0001
$ws.Range("M136").Value = -6303
$ws.Range("N136").Value = -15782.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 93051096
$ws.Range("I4").Value = 89008380
$ws.Range("J4").Value = 105179240
$ws.Range("K4").Value = 267025140
$ws.Range("L4").Value = 315537720
$ws.Range("M4").Value = -267025028
$ws.Range("N4").Value = -315537944

$ws.Range("H56").Value = 6196
$ws.Range("I56").Value = 6196
$ws.Range("K56").Value = 6196
$ws.Range("M56").Value = -5666

$ws.Range("H97").Value = 120100.4
$ws.Range("I97").Value = 150050.5
$ws.Range("K97").Value = 450151.5
$ws.Range("M97").Value = -449655.5

$ws.Range("H131").Value = 1564.1809
$ws.Range("J131").Value = 1577.069
$ws.Range("L131").Value = 4731.207
$ws.Range("N131").Value = -14811.207

$ws.Range("H139").Value = 1113483.9
$ws.Range("J139").Value = 3999.9167
$ws.Range("L139").Value = 11999.7501
$ws.Range("N139").Value = -22279.7501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 20825908
$ws.Range("I11").Value = 25228888
$ws.Range("J11").Value = 1012500
$ws.Range("K11").Value = 25228888
$ws.Range("L11").Value = 1012500
$ws.Range("M11").Value = -25228749
$ws.Range("N11").Value = -1012778

$ws.Range("H132").Value = 2037.836
$ws.Range("I132").Value = 1933
$ws.Range("K132").Value = 5799
$ws.Range("M132").Value = -3269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3671239
$ws.Range("I46").Value = 899.8
$ws.Range("J46").Value = 4982074.5
$ws.Range("K46").Value = 899.8
$ws.Range("L46").Value = 4982074.5
$ws.Range("M46").Value = -711.8
$ws.Range("N46").Value = -4982450.5

$ws.Range("H122").Value = 5694.625
$ws.Range("J122").Value = 2750
$ws.Range("L122").Value = 8250
$ws.Range("N122").Value = -13150

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws.Range("H136").Value = 5026.5835
$ws.Range("I136").Value = 3564.3
$ws.Range("J136").Value = 6071.0713
$ws.Range("K136").Value = 10692.9
$ws.Range("L136").Value = 18213.2139
$ws.Range("M136").Value = -8142.900000000001
$ws.Range("N136").Value = -23313.2139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1795.4546
$ws.Range("I96").Value = 826
$ws.Range("J96").Value = 2958.8
$ws.Range("K96").Value = 826
$ws.Range("L96").Value = 2958.8
$ws.Range("M96").Value = 547
$ws.Range("N96").Value = -5704.8

$ws.Range("H100").Value = 33168.332
$ws.Range("I100").Value = 21198.154
$ws.Range("K100").Value = 42396.308
$ws.Range("M100").Value = -41855.308

$ws.Range("H122").Value = 4188.3706
$ws.Range("I122").Value = 1748.9667
$ws.Range("K122").Value = 5246.9001
$ws.Range("M122").Value = -2796.9001

$ws.Range("H132").Value = 9489.870000000001
$ws.Range("I132").Value = 11454.641
$ws.Range("K132").Value = 34363.923
$ws.Range("M132").Value = -31833.923
